$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newUrl = "https://www.360dx.com/regulatory-news-fda-approvals/roche-abbott-inflammatix-paige-others-gain-510k-clearances-january"
$newKeyword = "digital pathology"
$newTitle = "Roche, Abbott, Inflammatix, Paige, Others Gain 510(k) Clearances in January"

$row = 16

$ws.Cells.Item($row, 2).Value = $newKeyword
$ws.Cells.Item($row, 3).Value = $newTitle

$aCell = $ws.Cells.Item($row, 1)
$aCell.Value = $newUrl
$ws.Hyperlinks.Add($aCell, $newUrl)
$aCell.Style = $ws.Cells.Item($row - 1, 1).Style
